$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update address for 영월군 row: "중앙1로 59" -> "영월군 중앙1로 59"
$ws.Range("C2").Value = "영월군 중앙1로 59"

# Update address for 하늘빛정신건강의학과의원 row: "원일로 150-1" -> "원주시 원일로 150-1"
$ws.Range("C15").Value = "원주시 원일로 150-1"

# Match final selection shown in the saved workbook
$ws.Range("C16").Select()
